$d = $word.ActiveDocument

# Locate the "Requisitos" bullet-list paragraph (the one holding the LOB/LOQ/LOM requirement codes)
$p = $null
foreach ($cand in $d.Paragraphs) {
    if ($cand.Range.Text.StartsWith("LOB1003")) {
        $p = $cand
        break
    }
}

$pStart = $p.Range.Start
$full = $p.Range.Text
$oldContentLen = $full.Length - 1   # exclude the trailing paragraph mark

# Target (reordered) lines, listed last-to-first: repeated InsertBefore calls on a
# range collapsed at pStart each push the new text to the front, so inserting in
# reverse yields the correct final top-to-bottom order.
$newTextReversed = @(
    "LOB1019 -  Física II  (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOB1041 -  Física Experimental II  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1024 -  Mecânica  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1012 -  Estatística  (Requisito)",
    "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1037 -  Àlgebra Linear  (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOQ4251 -  Fundamentos de Química  (Requisito)",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)"
)

# Insert each reordered requirement line (plus its line break) as its own fresh run
# immediately before the old list - this keeps every line in a distinct <w:r>.
$insertPoint = $d.Range($pStart, $pStart)
foreach ($line in $newTextReversed) {
    $insertPoint.InsertBefore($line + [char]11)
}

# The old (pre-reorder) list text now sits right after all our insertions
$oldStart = $insertPoint.End
$oldRange = $d.Range($oldStart, $oldStart + $oldContentLen)
$oldRange.Delete()

Write-Output "reordered requisitos list"
